# Session 1 - Hooking Up: update the IDE's bullet on the Python setup slide
# and refresh the "datetimeFigureOut" date placeholders on the slide master
# and every slide layout (PowerPoint re-stamps these automatically whenever
# the deck is re-saved on a different day).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 ("Setup - Python"): "Many open-source IDE's" -> split into
#    "Many " / "free" / " " / "IDE's and " / "Text Editors"
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(6)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(5, 1)

# "Many open-source IDE's"
#  01234567890123456789012
# Replace "open-source" (chars 6-16, 1-based Start+5 len 11) with "free"
$mid = $tr.Characters($para.Start + 5, 11)
$mid.Text = "free"

# Re-fetch paragraph (text shifted) and split the joining space into its own run
$para = $tr.Paragraphs(5, 1)
$spaceRun = $tr.Characters($para.Start + 9, 1)
$spaceRun.Text = " "

# Re-fetch paragraph; extend "IDE's" in place to "IDE's and " (keeps the
# original run's formatting -- this is the only run the diff leaves with
# dirty="0")
$para = $tr.Paragraphs(5, 1)
$ideRun = $tr.Characters($para.Start + 10, 5)
$ideRun.Text = "IDE" + [char]0x2019 + "s and "

# Append the new trailing run "Text Editors"
$para = $tr.Paragraphs(5, 1)
$para.InsertAfter("Text Editors")

# ---------------------------------------------------------------------------
# 2) Refresh the date placeholder text on the slide master + every layout
# ---------------------------------------------------------------------------
$newDate = "3/5/2019"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -eq "Date Placeholder 3") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -eq "Date Placeholder 3") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
